$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: B 14->16, C 15->16 ---
# Excel's COM ColumnWidth property is offset from the raw OOXML <col width>
# value by the default-font padding constant (~0.8333 chars), so subtract it
# to land on an exact raw width of 16.
$rawToComOffset = 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth = 16 - $rawToComOffset
$ws.Columns.Item(3).ColumnWidth = 16 - $rawToComOffset

# --- Keep the revenue-table numbers as plain text (matches the source
#     workbook, which stores these as text, not numeric, cells) ---
$ws.Range("B4:C7").NumberFormat = "@"

# --- Capture/propagate formatting before touching values ---
# A8 ("Total Sale of Services") currently carries the bold+bordered style.
# A4 and A7 need to become bold+bordered in the final layout ("Particulars"
# header row and the new "Total Sale of Services" row respectively), so copy
# that format onto them first (A5/A6 already carry the plain bordered style
# that the "Domestic"/"Exports" rows need, so they are left untouched).
$ws.Range("A8").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 4: new "Particulars" header row ---
$ws.Range("A4").Value = "Particulars"
$ws.Range("B4").Value = "March 31, 2024"
$ws.Range("C4").Value = "March 31, 2023"

# --- Row 5: Domestic ---
$ws.Range("A5").Value = "Domestic"
$ws.Range("B5").Value = "354.42"
$ws.Range("C5").Value = "0.00"

# --- Row 6: Exports ---
$ws.Range("A6").Value = "Exports"
$ws.Range("B6").Value = "10.00"
$ws.Range("C6").Value = "0.00"

# --- Row 7: Total Sale of Services ---
$ws.Range("A7").Value = "Total Sale of Services"
$ws.Range("B7").Value = "364.42"
$ws.Range("C7").Value = "0.00"

# --- Old row 8 ("Total Sale of Services") is no longer needed; clear it
#     completely (contents + formatting) so it becomes a blank spacer row ---
$ws.Range("A8:C8").Clear()

# --- Row 9 was a blank spacer row; deleting it shifts rows 10-13
#     ("Summary:", "Total Amount:", "Total Amount (Lakhs):",
#     "Matched Accounts Count:") up into rows 9-12, which matches the
#     target layout (and shrinks the used range to A1:C12). ---
$ws.Rows.Item(9).Delete()
